# "Demo nach ganz hinten" - move the "Demo" slide to the very end of the
# deck (just before the final "Link" slide).
$p = $ppt.ActivePresentation

# The "Demo" slide is currently slide 13; the deck has 21 slides total and
# the very last slide ("Link") should stay last, so "Demo" needs to land
# at position 20 (i.e. directly before the current last slide).
$demoSlide = $p.Slides.Item(13)
$demoSlide.MoveTo($p.Slides.Count - 1)
